# Generate Report for Handback
#
# The e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md file has now been handed back
# (for both zh-cn and de-de). Update the Overview sheet status for that row,
# and fill in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: rows for e7b04fa2-...md (row 3) and f08fbba8-...md (row 4)
# move from "Ready for handoff" to "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C4").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 3 (e7b04fa2-...)
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("F3").Value = "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f342d0d35ce2b9e7fe029d7c4f2cb041ba114462/e2e/e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md", "", "", "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md") | Out-Null
$wsZh.Range("G3").Value = "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f342d0d35ce2b9e7fe029d7c4f2cb041ba114462/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.zh-cn.xlf", "", "", "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.zh-cn.xlf") | Out-Null
$wsZh.Range("H3").Value = "2016-03-14 09:22:03"

# Row 4 (f08fbba8-..., shares the same target/handback file as row 3)
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("F4").Value = "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f342d0d35ce2b9e7fe029d7c4f2cb041ba114462/e2e/e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md", "", "", "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md") | Out-Null
$wsZh.Range("G4").Value = "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f342d0d35ce2b9e7fe029d7c4f2cb041ba114462/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.zh-cn.xlf", "", "", "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.zh-cn.xlf") | Out-Null
$wsZh.Range("H4").Value = "2016-03-14 09:22:03"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 3 (e7b04fa2-...)
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("F3").Value = "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a4b1d3cdfeeb8dbdebcbe3fe72b728c9270d75b0/e2e/e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md", "", "", "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md") | Out-Null
$wsDe.Range("G3").Value = "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a4b1d3cdfeeb8dbdebcbe3fe72b728c9270d75b0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.de-de.xlf", "", "", "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.de-de.xlf") | Out-Null
$wsDe.Range("H3").Value = "2016-03-14 09:22:16"

# Row 4 (f08fbba8-..., shares the same target/handback file as row 3)
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("F4").Value = "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a4b1d3cdfeeb8dbdebcbe3fe72b728c9270d75b0/e2e/e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md", "", "", "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.md") | Out-Null
$wsDe.Range("G4").Value = "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a4b1d3cdfeeb8dbdebcbe3fe72b728c9270d75b0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.de-de.xlf", "", "", "e7b04fa2-7ac6-4e7e-b65e-2adfe8b500c3.8059e2223553327a2ecc870e35961e0ea1b06037.de-de.xlf") | Out-Null
$wsDe.Range("H4").Value = "2016-03-14 09:22:16"

Write-Host "Handback report generated."
